# Bitacora_tareas.xlsx -> v1.6 update
# - Log: add rows 33-37
# - Resumen: add rows 35-39
# - Versiones: add row 8
# - Presupuesto: insert two new rubros (rows 4-5), renumbering existing rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Log" - add rows 33..37
# ---------------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

$log.Range("A33").Value = "27/02/2025"
$log.Range("B33").Value = "19:00"
$log.Range("C33").Value = "Regla bitácora: actualizar todas las solapas necesarias"
$log.Range("D33").Value = "La regla pasa a exigir actualizar todas las solapas que correspondan: Log, Resumen (si aplica), Presupuesto (cuando la tarea agrega o cambia un entregable comercial), Versiones (en despliegue). Presupuesto se actualiza con el rubro ""Detección de duplicados y gestión de errores""."
$log.Range("E33").Value = "Diagnostico"

$log.Range("A34").Value = "27/02/2025"
$log.Range("B34").Value = "19:10"
$log.Range("C34").Value = "Solapa Evolución (tabla dinámica)"
$log.Range("D34").Value = "Nueva pestaña Evolución: tabla dinámica con Agrupar por (Categoría o Cuenta contable) como fila y Período (Diario o Mensual) como columna. Diario muestra fecha (día), Mensual muestra MM-YYYY. Celdas = neto (ingresos - egresos) en la moneda seleccionada. Columna Total por fila."
$log.Range("E34").Value = "Diagnostico"

$log.Range("A35").Value = "27/02/2025"
$log.Range("B35").Value = "19:20"
$log.Range("C35").Value = "Evolución: clic en valor y exportar a Excel"
$log.Range("D35").Value = "Al hacer clic en un valor de la tabla Evolución se abre un modal con detalle mínimo: Fecha, Categoría, Descripción, Monto (registros que componen esa celda). Botón Exportar Evolución a Excel exporta la tabla resultante según los filtros Agrupar por y Período."
$log.Range("E35").Value = "Diagnostico"

$log.Range("A36").Value = "27/02/2025"
$log.Range("B36").Value = "19:30"
$log.Range("C36").Value = "Exportaciones: título moneda, icono Excel, Exportar Base Histórica"
$log.Range("D36").Value = "En todas las exportaciones a Excel se agrega una fila título que indica la moneda (o que ver columna moneda). Icono tipo Excel (tabla/grid) en botones de exportar. Exportar base de transacciones movido a la línea del selector de moneda con título ""Exportar Base Histórica"" e icono Excel; mismo icono en Exportar Evolución a Excel."
$log.Range("E36").Value = "Diagnostico"

$log.Range("A37").Value = "27/02/2025"
$log.Range("B37").Value = "19:40"
$log.Range("C37").Value = "Evolución: ingreso primero, luego egreso"
$log.Range("D37").Value = "En la tabla Evolución las filas (categorías o cuentas) se ordenan primero las de ingreso (total >= 0) y luego las de egreso (total < 0); dentro de cada grupo orden alfabético."
$log.Range("E37").Value = "Diagnostico"

# ---------------------------------------------------------------------------
# Sheet "Resumen" - add rows 35..39
# ---------------------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")

$resumen.Range("A35").Value = "Regla bitácora"
$resumen.Range("B35").Value = "Actualizar todas las solapas necesarias: Log (siempre que haya tarea), Resumen (si cambia funcionalidad), Presupuesto (si agrega o cambia entregable comercial), Versiones (en despliegue). Regenerar Excel tras editar crear-bitacora-excel.js."

$resumen.Range("A36").Value = "Evolución (tabla dinámica)"
$resumen.Range("B36").Value = "Solapa Evolución: Agrupar por = Categoría o Cuenta contable (fila); Período = Diario (fecha por día) o Mensual (MM-YYYY). Columnas = períodos, celdas = neto en moneda seleccionada, columna Total."

$resumen.Range("A37").Value = "Evolución: detalle al clic y exportar"
$resumen.Range("B37").Value = "Clic en un valor de la tabla Evolución abre modal con detalle: Fecha, Categoría, Descripción, Monto. Exportar Evolución a Excel exporta la tabla según filtros Agrupar por y Período."

$resumen.Range("A38").Value = "Exportaciones Excel"
$resumen.Range("B38").Value = "Todas las exportaciones incluyen una fila título con la moneda. Exportar Base Histórica (icono Excel) en la línea del selector de moneda; Exportar Evolución a Excel con el mismo icono."

$resumen.Range("A39").Value = "Evolución: orden ingreso/egreso"
$resumen.Range("B39").Value = "En la tabla Evolución las filas se muestran primero las de ingreso (total >= 0) y luego las de egreso (total < 0); dentro de cada grupo orden alfabético. Aplica tanto al agrupar por Categoría como por Cuenta contable."

# ---------------------------------------------------------------------------
# Sheet "Versiones" - add row 8
# ---------------------------------------------------------------------------
$versiones = $wb.Worksheets.Item("Versiones")

# Force text storage so "1.6" isn't auto-converted to a number (matches A2:A7)
$versiones.Range("A8").NumberFormat = "@"
$versiones.Range("A8").Value = "1.6"
$versiones.Range("B8").Value = "27/02/2025"
$versiones.Range("C8").Value = "Export Excel: botones verde y blanco; Evolución: orden ingreso luego egreso; modal detalle Evolución con columna Origen y modal más ancho"

# ---------------------------------------------------------------------------
# Sheet "Presupuesto" - insert two new rubros before "Bitácora y documentación"
# ---------------------------------------------------------------------------
$presupuesto = $wb.Worksheets.Item("Presupuesto")

# Insert two blank rows at position 4 (existing rows 4-6 shift to 6-8)
$presupuesto.Rows.Item(4).Insert()
$presupuesto.Rows.Item(4).Insert()
Write-Host "Presupuesto: filas insertadas"

$presupuesto.Range("A4").Value = "Detección de duplicados y gestión de errores"
$presupuesto.Range("B4").Value = "Detección de potencial duplicado (fecha, monto, tipo, cliente, descripción similar), tipo de error (inconsistencia / duplicado), filtro por tipo, modal de comparación con id_origen y Cliente, acciones anular o eliminar registro."
$presupuesto.Range("C4").Value = 85000

$presupuesto.Range("A5").Value = "Evolución (tabla dinámica)"
$presupuesto.Range("B5").Value = "Solapa Evolución: tabla dinámica con filas por Categoría o Cuenta contable y columnas por Período (Diario o Mensual). Neto por celda en moneda seleccionada."
$presupuesto.Range("C5").Value = 55000

Write-Host "Edit complete: Log, Resumen, Versiones, Presupuesto actualizados"
